$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.926.71"
$ws.Range("E2").Value = "  -3.77%  "
$ws.Range("D3").Value = "3.467.30"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'603.29"
$ws.Range("E5").Value = "  -4.29%  "
$ws.Range("D6").Value = "'147.76"
$ws.Range("E6").Value = "  -6.69%  "
$ws.Range("D7").Value = "3.465.08"
$ws.Range("E7").Value = "  -4.00%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.484"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("E10").Value = "  -4.75%  "
$ws.Range("D11").Value = "'7.48"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  -4.14%  "
$ws.Range("D13").Value = "'0.0000214"
$ws.Range("E13").Value = "  -5.95%  "
$ws.Range("D14").Value = "'31.67"
$ws.Range("D15").Value = "4.060.19"
$ws.Range("E15").Value = "  -3.95%  "
$ws.Range("D16").Value = "3.473.43"
$ws.Range("E16").Value = "  -3.76%  "
$ws.Range("D17").Value = "66.934.48"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  -4.59%  "
$ws.Range("D20").Value = "'15.32"
$ws.Range("E20").Value = "  -5.11%  "
$ws.Range("D21").Value = "'10.02"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "'439.93"
$ws.Range("E22").Value = "  -4.74%  "
$ws.Range("D23").Value = "'0.608"
$ws.Range("E23").Value = "  -5.78%  "
$ws.Range("D24").Value = "'78.60"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "3.610.11"
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("E27").Value = "  -9.89%  "
$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  -7.92%  "
$ws.Range("D29").Value = "'8.39"
$ws.Range("E29").Value = "  -9.70%  "
$ws.Range("E30").Value = "  -6.56%  "
$ws.Range("D31").Value = "'1.61"
$ws.Range("E31").Value = "  -6.56%  "
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").Value = "'6.08"
$ws.Range("E35").Value = "  -7.17%  "
$ws.Range("D36").Value = "3.465.57"
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  -7.64%  "
$ws.Range("E38").Value = "  -6.64%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.18"
$ws.Range("E41").Value = "  -9.66%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'172.38"
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("D43").Value = "'0.0888"
$ws.Range("E43").Value = "  -4.16%  "
$ws.Range("D44").Value = "'5.39"
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("D45").Value = "'0.884"
$ws.Range("E45").Value = "  -3.22%  "
$ws.Range("D46").Value = "'29.20"
$ws.Range("E46").Value = "  -9.11%  "
$ws.Range("D47").Value = "'46.15"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  -11.00%  "
$ws.Range("D49").Value = "'7.47"
$ws.Range("E49").Value = "  -4.55%  "
$ws.Range("E50").Value = "  -10.60%  "
$ws.Range("D51").Value = "'0.986"
$ws.Range("E51").Value = "  -4.93%  "
